$d = $word.ActiveDocument

# --- Locate the paragraph ending "...Dir() ... newly created class." ---
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Dir()*newly created class.*") {
        $targetIndex = $i
        break
    }
}

# --- Append a new run to the end of that paragraph (before its paragraph mark) ---
$tp = $d.Paragraphs.Item($targetIndex)
$tr = $tp.Range
$work = $d.Range($tr.Start, $tr.End - 1)
$work.Collapse(0)
$work.InsertAfter(' Object creation – constructor and destruction – destructor.')
$work.Font.Reset()

$curIndex = $targetIndex

# --- New paragraph 1 ---
$cp = $d.Paragraphs.Item($curIndex)
$cr = $cp.Range
$cw = $d.Range($cr.Start, $cr.End - 1)
$cw.Collapse(0)
$cw.InsertParagraphAfter()
$curIndex = $curIndex + 1
$np = $d.Paragraphs.Item($curIndex)
$nr = $np.Range
$nw = $d.Range($nr.Start, $nr.End - 1)
$nw.InsertAfter('When an object is created it runs all the variables and then runs the constructor ie……………..  ')
$nw.Font.Reset()
$nw.Collapse(0)
$nw.InsertAfter('def __init__(self):     self.x = 0')
$nw.Font.Reset()
$nw.Collapse(0)
$nw.InsertAfter(' ')
$nw.Font.Reset()

# --- New paragraph 2 ---
$cp = $d.Paragraphs.Item($curIndex)
$cr = $cp.Range
$cw = $d.Range($cr.Start, $cr.End - 1)
$cw.Collapse(0)
$cw.InsertParagraphAfter()
$curIndex = $curIndex + 1
$np = $d.Paragraphs.Item($curIndex)
$nr = $np.Range
$nw = $d.Range($nr.Start, $nr.End - 1)
$nw.InsertAfter('Eg.  ')
$nw.Font.Reset()
$nw.Collapse(0)
$nw.InsertAfter('an = PartyAnimal')
$nw.Font.Reset()
$nw.Collapse(0)
$nw.InsertAfter('()')
$nw.Font.Reset()
$nw.Collapse(0)
$nw.InsertAfter(' I want to involved when this object is created')
$nw.Font.Color = 255
$nw.Collapse(0)
$nw.InsertAfter(' – this executes till I am constructed message. Then when created object loses its pointer somehow. Ie by maybe assigning it a value ie. an = 20. Then destructor gets activated ')
$nw.Font.Reset()
$nw.Collapse(0)
$nw.InsertAfter('(')
$nw.Font.Reset()
$nw.Collapse(0)
$nw.InsertAfter('I want to be involved when the object is destroyed')
$nw.Font.Color = 255
$nw.Collapse(0)
$nw.InsertAfter(') ')
$nw.Font.Reset()
$nw.Collapse(0)
$nw.InsertAfter('ie ')
$nw.Font.Reset()
$nw.Collapse(0)
$nw.InsertAfter('  def __del__(self):      print(''I am destructed'', self.x)')
$nw.Font.Reset()

# --- New paragraph 3 ---
$cp = $d.Paragraphs.Item($curIndex)
$cr = $cp.Range
$cw = $d.Range($cr.Start, $cr.End - 1)
$cw.Collapse(0)
$cw.InsertParagraphAfter()
$curIndex = $curIndex + 1
$np = $d.Paragraphs.Item($curIndex)
$nr = $np.Range
$nw = $d.Range($nr.Start, $nr.End - 1)
$nw.InsertAfter('Constructors can have additional parameters that can be used to set up instance variables for the particular instance of the class.')
$nw.Font.Reset()

# --- New paragraph 4 ---
$cp = $d.Paragraphs.Item($curIndex)
$cr = $cp.Range
$cw = $d.Range($cr.Start, $cr.End - 1)
$cw.Collapse(0)
$cw.InsertParagraphAfter()
$curIndex = $curIndex + 1
$np = $d.Paragraphs.Item($curIndex)
$nr = $np.Range
$nw = $d.Range($nr.Start, $nr.End - 1)
$nw.InsertAfter('Inheritance – when we make a new class we reuse the existing class and inherit all the capabilities of an existing class and then add our own little bit to make our new class. Basically called as “subclassing”, were attributes are inherited from their parent class and introduced their own.')
$nw.Font.Reset()

# --- New paragraph 5 ---
$cp = $d.Paragraphs.Item($curIndex)
$cr = $cp.Range
$cw = $d.Range($cr.Start, $cr.End - 1)
$cw.Collapse(0)
$cw.InsertParagraphAfter()
$curIndex = $curIndex + 1
$np = $d.Paragraphs.Item($curIndex)
$nr = $np.Range
$nw = $d.Range($nr.Start, $nr.End - 1)
$nw.InsertAfter('So basically inheritance will take the constructor ok parent class and also execute the constructor of current class. Same goes for destructor.')
$nw.Font.Reset()

# --- New paragraph 6 ---
$cp = $d.Paragraphs.Item($curIndex)
$cr = $cp.Range
$cw = $d.Range($cr.Start, $cr.End - 1)
$cw.Collapse(0)
$cw.InsertParagraphAfter()
$curIndex = $curIndex + 1
$np = $d.Paragraphs.Item($curIndex)
$nr = $np.Range
$nw = $d.Range($nr.Start, $nr.End - 1)
$nw.InsertAfter("X")
$nw.Delete()

# --- New paragraph 7 ---
$cp = $d.Paragraphs.Item($curIndex)
$cr = $cp.Range
$cw = $d.Range($cr.Start, $cr.End - 1)
$cw.Collapse(0)
$cw.InsertParagraphAfter()
$curIndex = $curIndex + 1
$np = $d.Paragraphs.Item($curIndex)
$nr = $np.Range
$nw = $d.Range($nr.Start, $nr.End - 1)
$nw.InsertAfter("X")
$nw.Delete()

# --- New paragraph 8 ---
$cp = $d.Paragraphs.Item($curIndex)
$cr = $cp.Range
$cw = $d.Range($cr.Start, $cr.End - 1)
$cw.Collapse(0)
$cw.InsertParagraphAfter()
$curIndex = $curIndex + 1
$np = $d.Paragraphs.Item($curIndex)
$nr = $np.Range
$nw = $d.Range($nr.Start, $nr.End - 1)
$nw.InsertAfter("X")
$nw.Delete()

# --- New paragraph 9 ---
$cp = $d.Paragraphs.Item($curIndex)
$cr = $cp.Range
$cw = $d.Range($cr.Start, $cr.End - 1)
$cw.Collapse(0)
$cw.InsertParagraphAfter()
$curIndex = $curIndex + 1
$np = $d.Paragraphs.Item($curIndex)
$nr = $np.Range
$nw = $d.Range($nr.Start, $nr.End - 1)
$nw.InsertAfter("X")
$nw.Delete()
